$wb = $excel.ActiveWorkbook

$dbd = $wb.Worksheets.Item("DBD")
$dbs = $wb.Worksheets.Item("DBS")

# Add the new "findSupNoEntdy" lookup row to the DBS sheet, copying the
# formatting of the row above (row 3) so fonts/borders/alignment match.
$dbs.Range("A3:C3").Copy()
$dbs.Range("A5:C5").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

$dbs.Range("A5").Value = "findSupNoEntdy"
$dbs.Range("B5").Value = "Entdy >= ,AND Entdy <=,AND SupNo %"
$dbs.Range("C5").Value = "SupNo ASC,Entdy ASC"

# Update selections to match the saved state.
$dbd.Range("B14").Select() | Out-Null
$dbs.Range("B14").Select() | Out-Null

# The DBS sheet becomes the active tab on save.
$dbs.Activate() | Out-Null
